$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Import" row (row 2) now needs to list both POJO classes used by the
# rule sheet, not just Policy.
$ws.Range("B2").Value = "com.redhat.prudential_poc.pojo.Insured,com.redhat.prudential_poc.pojo.Policy"

# Move/restore the active selection to B3 (matches the saved cursor position).
$ws.Range("B3").Select() | Out-Null
